$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 286.09
$ws.Cells.Item(15, 9).Value = 286.09
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 858.27
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = -689.27

$ws.Cells.Item(17, 8).Value = 531725.6
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 531725.6
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 1595176.8
$ws.Cells.Item(17, 14).Value = -1595512.8

$ws.Cells.Item(19, 8).Value = 201086.94
$ws.Cells.Item(19, 9).Value = 387980.4
$ws.Cells.Item(19, 10).Value = 843.9286
$ws.Cells.Item(19, 11).Value = 387980.4
$ws.Cells.Item(19, 12).Value = 843.9286
$ws.Cells.Item(19, 13).Value = -387805.4
$ws.Cells.Item(19, 14).Value = -1193.9286

$ws.Cells.Item(28, 8).Value = 2425.3333
$ws.Cells.Item(28, 9).Value = 3168.75
$ws.Cells.Item(28, 10).Value = 938.5
$ws.Cells.Item(28, 11).Value = 3168.75
$ws.Cells.Item(28, 12).Value = 938.5
$ws.Cells.Item(28, 13).Value = -2683.75
$ws.Cells.Item(28, 14).Value = -1908.5

$ws.Cells.Item(92, 8).Value = 1176.95
$ws.Cells.Item(92, 9).Value = 1411.5
$ws.Cells.Item(92, 10).Value = 238.75
$ws.Cells.Item(92, 11).Value = 1411.5
$ws.Cells.Item(92, 12).Value = 238.75
$ws.Cells.Item(92, 13).Value = -163.5
$ws.Cells.Item(92, 14).Value = -2734.75

$ws.Cells.Item(107, 8).Value = 521.5
$ws.Cells.Item(107, 9).Value = 528.7778
$ws.Cells.Item(107, 10).Value = 456
$ws.Cells.Item(107, 11).Value = 528.7778
$ws.Cells.Item(107, 12).Value = 456
$ws.Cells.Item(107, 13).Value = 1391.2222
$ws.Cells.Item(107, 14).Value = -4296

$ws.Cells.Item(112, 8).Value = 1619.375
$ws.Cells.Item(112, 9).Value = 1036.6666
$ws.Cells.Item(112, 10).Value = 1753.8462
$ws.Cells.Item(112, 11).Value = 3109.9998
$ws.Cells.Item(112, 12).Value = 5261.5386
$ws.Cells.Item(112, 13).Value = -2001.9998
$ws.Cells.Item(112, 14).Value = -7477.5386

$ws.Cells.Item(116, 8).Value = 3628.2693
$ws.Cells.Item(116, 9).Value = 3451.7058
$ws.Cells.Item(116, 10).Value = 3961.7778
$ws.Cells.Item(116, 11).Value = 3451.7058
$ws.Cells.Item(116, 12).Value = 3961.7778
$ws.Cells.Item(116, 13).Value = -9.705800000000181
$ws.Cells.Item(116, 14).Value = -10845.7778

$ws.Cells.Item(129, 8).Value = 918330.25
$ws.Cells.Item(129, 9).Value = 274.66666
$ws.Cells.Item(129, 10).Value = 1262601.1
$ws.Cells.Item(129, 11).Value = 823.9999799999999
$ws.Cells.Item(129, 12).Value = 3787803.3
$ws.Cells.Item(129, 13).Value = 4176.00002
$ws.Cells.Item(129, 14).Value = -3797803.3

$ws.Cells.Item(137, 8).Value = 646.73334
$ws.Cells.Item(137, 9).Value = 646.2308
$ws.Cells.Item(137, 10).Value = 650
$ws.Cells.Item(137, 11).Value = 1938.6924
$ws.Cells.Item(137, 12).Value = 1950
$ws.Cells.Item(137, 13).Value = 611.3075999999999
$ws.Cells.Item(137, 14).Value = -7050

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2708.4
$ws.Cells.Item(32, 9).Value = 2584.2424
$ws.Cells.Item(32, 10).Value = 15000
$ws.Cells.Item(32, 11).Value = 2584.2424
$ws.Cells.Item(32, 12).Value = 15000
$ws.Cells.Item(32, 13).Value = -2297.2424

$ws.Cells.Item(61, 8).Value = 2029.2593
$ws.Cells.Item(61, 9).Value = 1991.6666
$ws.Cells.Item(61, 10).Value = 2330
$ws.Cells.Item(61, 11).Value = 1991.6666
$ws.Cells.Item(61, 12).Value = 2330
$ws.Cells.Item(61, 13).Value = -1779.6666
$ws.Cells.Item(61, 14).Value = -2754

$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).ClearContents()

$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 2029.2593
$ws.Cells.Item(136, 9).Value = 1991.6666
$ws.Cells.Item(136, 10).Value = 2330
$ws.Cells.Item(136, 11).Value = 5974.9998
$ws.Cells.Item(136, 12).Value = 6990
$ws.Cells.Item(136, 13).Value = -3424.9998
$ws.Cells.Item(136, 14).Value = -12090

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 25644118
$ws.Cells.Item(86, 9).Value = 29414592
$ws.Cells.Item(86, 10).Value = 4902
$ws.Cells.Item(86, 11).Value = 29414592
$ws.Cells.Item(86, 12).Value = 4902
$ws.Cells.Item(86, 13).Value = -29413469
$ws.Cells.Item(86, 14).Value = -7148

$ws.Cells.Item(89, 8).Value = 25644118
$ws.Cells.Item(89, 9).Value = 29414592
$ws.Cells.Item(89, 10).Value = 4902
$ws.Cells.Item(89, 11).Value = 147072960
$ws.Cells.Item(89, 12).Value = 24510
$ws.Cells.Item(89, 13).Value = -147067344
$ws.Cells.Item(89, 14).Value = -35742

$ws.Cells.Item(99, 8).Value = 66668720
$ws.Cells.Item(99, 9).Value = 83335320
$ws.Cells.Item(99, 10).Value = 2333.3333
$ws.Cells.Item(99, 11).Value = 83335320
$ws.Cells.Item(99, 12).Value = 2333.3333
$ws.Cells.Item(99, 13).Value = -83333822
$ws.Cells.Item(99, 14).Value = -5329.3333

$ws.Cells.Item(134, 8).Value = 27259.25
$ws.Cells.Item(134, 9).Value = 2077.923
$ws.Cells.Item(134, 10).Value = 74024.57000000001
$ws.Cells.Item(134, 11).Value = 6233.768999999999
$ws.Cells.Item(134, 12).Value = 222073.71
$ws.Cells.Item(134, 13).Value = -3698.768999999999
$ws.Cells.Item(134, 14).Value = -227143.71

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1473.3
$ws.Cells.Item(16, 9).Value = 1205.5
$ws.Cells.Item(16, 10).Value = 1875
$ws.Cells.Item(16, 11).Value = 1205.5
$ws.Cells.Item(16, 12).Value = 1875
$ws.Cells.Item(16, 13).Value = -918.5
$ws.Cells.Item(16, 14).Value = -2449

$ws.Cells.Item(31, 8).Value = 3154.0833
$ws.Cells.Item(31, 9).Value = 3090.3809
$ws.Cells.Item(31, 10).Value = 3600
$ws.Cells.Item(31, 11).Value = 3090.3809
$ws.Cells.Item(31, 12).Value = 3600
$ws.Cells.Item(31, 13).Value = -2795.3809
$ws.Cells.Item(31, 14).Value = -4190

$ws.Cells.Item(34, 8).Value = 3154.0833
$ws.Cells.Item(34, 9).Value = 3090.3809
$ws.Cells.Item(34, 10).Value = 3600
$ws.Cells.Item(34, 11).Value = 3090.3809
$ws.Cells.Item(34, 12).Value = 3600
$ws.Cells.Item(34, 13).Value = -2888.3809
$ws.Cells.Item(34, 14).Value = -4004

$ws.Cells.Item(68, 8).Value = 15777.223
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 15777.223
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 15777.223
$ws.Cells.Item(68, 14).Value = -17275.223

$ws.Cells.Item(71, 8).Value = 15777.223
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 15777.223
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 47331.669
$ws.Cells.Item(71, 14).Value = -54819.669

$ws.Cells.Item(113, 8).Value = 1473.3
$ws.Cells.Item(113, 9).Value = 1205.5
$ws.Cells.Item(113, 10).Value = 1875
$ws.Cells.Item(113, 11).Value = 1205.5
$ws.Cells.Item(113, 12).Value = 1875
$ws.Cells.Item(113, 13).Value = 964.5
$ws.Cells.Item(113, 14).Value = -6215

$ws.Cells.Item(134, 8).Value = 14707023
$ws.Cells.Item(134, 9).Value = 1007.0417
$ws.Cells.Item(134, 10).Value = 50001460
$ws.Cells.Item(134, 11).Value = 3021.1251
$ws.Cells.Item(134, 12).Value = 150004380
$ws.Cells.Item(134, 13).Value = -486.1251000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 25516.75
$ws.Cells.Item(92, 9).Value = 100434
$ws.Cells.Item(92, 10).Value = 544.3333
$ws.Cells.Item(92, 11).Value = 301302
$ws.Cells.Item(92, 12).Value = 1632.9999
$ws.Cells.Item(92, 13).Value = -300054
$ws.Cells.Item(92, 14).Value = -4128.9999

$ws.Cells.Item(136, 8).Value = 94832.63
$ws.Cells.Item(136, 9).Value = 501499.5
$ws.Cells.Item(136, 10).Value = 4462.222
$ws.Cells.Item(136, 11).Value = 1504498.5
$ws.Cells.Item(136, 12).Value = 13386.666
$ws.Cells.Item(136, 13).Value = -1499398.5
$ws.Cells.Item(136, 14).Value = -23586.666

$ws.Cells.Item(138, 8).Value = 3762.8572
$ws.Cells.Item(138, 9).Value = 2120
$ws.Cells.Item(138, 10).Value = 4995
$ws.Cells.Item(138, 11).Value = 6360
$ws.Cells.Item(138, 12).Value = 14985
$ws.Cells.Item(138, 13).Value = -1220

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 6695002.5
$ws.Cells.Item(24, 9).Value = 10000000
$ws.Cells.Item(24, 10).Value = 85007
$ws.Cells.Item(24, 11).Value = 10000000
$ws.Cells.Item(24, 12).Value = 85007
$ws.Cells.Item(24, 13).Value = -9999827
$ws.Cells.Item(24, 14).Value = -85353

$ws.Cells.Item(93, 8).Value = 12000
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 12000
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 12000
$ws.Cells.Item(93, 14).Value = -15744

$ws.Cells.Item(107, 8).Value = 376.54166
$ws.Cells.Item(107, 9).Value = 265.10526
$ws.Cells.Item(107, 10).Value = 800
$ws.Cells.Item(107, 11).Value = 265.10526
$ws.Cells.Item(107, 12).Value = 800
$ws.Cells.Item(107, 13).Value = 1654.89474
$ws.Cells.Item(107, 14).Value = -4640

$ws.Cells.Item(122, 8).Value = 1317484.4
$ws.Cells.Item(122, 9).Value = 1881299.1
$ws.Cells.Item(122, 10).Value = 1916.6666
$ws.Cells.Item(122, 11).Value = 5643897.300000001
$ws.Cells.Item(122, 12).Value = 5749.9998
$ws.Cells.Item(122, 13).Value = -5641447.300000001
$ws.Cells.Item(122, 14).Value = -10649.9998

$ws.Cells.Item(132, 8).Value = 2777.1628
$ws.Cells.Item(132, 9).Value = 2675.2334
$ws.Cells.Item(132, 10).Value = 3012.3845
$ws.Cells.Item(132, 11).Value = 8025.7002
$ws.Cells.Item(132, 12).Value = 9037.1535
$ws.Cells.Item(132, 13).Value = -5495.7002
$ws.Cells.Item(132, 14).Value = -14097.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1185.2858
$ws.Cells.Item(22, 9).Value = 400
$ws.Cells.Item(22, 10).Value = 1245.6923
$ws.Cells.Item(22, 11).Value = 400
$ws.Cells.Item(22, 12).Value = 1245.6923
$ws.Cells.Item(22, 13).Value = -105
$ws.Cells.Item(22, 14).Value = -1835.6923

$ws.Cells.Item(27, 8).Value = 1185.2858
$ws.Cells.Item(27, 9).Value = 400
$ws.Cells.Item(27, 10).Value = 1245.6923
$ws.Cells.Item(27, 11).Value = 400
$ws.Cells.Item(27, 12).Value = 1245.6923
$ws.Cells.Item(27, 13).Value = -293
$ws.Cells.Item(27, 14).Value = -1459.6923

$ws.Cells.Item(46, 8).Value = 50001700
$ws.Cells.Item(46, 9).Value = 50001700
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 50001700
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = -50001512
$ws.Cells.Item(46, 14).ClearContents()

$ws.Cells.Item(58, 8).Value = 1000
$ws.Cells.Item(58, 9).Value = 1000
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 1000
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -740

$ws.Cells.Item(61, 8).Value = 1329.375
$ws.Cells.Item(61, 9).Value = 1351.5385
$ws.Cells.Item(61, 10).Value = 1233.3334
$ws.Cells.Item(61, 11).Value = 1351.5385
$ws.Cells.Item(61, 12).Value = 1233.3334
$ws.Cells.Item(61, 13).Value = -1149.5385
$ws.Cells.Item(61, 14).Value = -1637.3334

$ws.Cells.Item(62, 8).Value = 19999.5
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 19999.5
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 19999.5
$ws.Cells.Item(62, 14).Value = -21247.5

$ws.Cells.Item(65, 8).Value = 19999.5
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 19999.5
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 59998.5
$ws.Cells.Item(65, 14).Value = -66238.5

$ws.Cells.Item(76, 8).Value = 10000
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 10000
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 10000
$ws.Cells.Item(76, 13).ClearContents()
$ws.Cells.Item(76, 14).Value = -10676

$ws.Cells.Item(79, 8).Value = 10000
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 10000
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 10000
$ws.Cells.Item(79, 13).ClearContents()
$ws.Cells.Item(79, 14).Value = -12340

$ws.Cells.Item(100, 8).Value = 7015624
$ws.Cells.Item(100, 9).Value = 7483245.5
$ws.Cells.Item(100, 10).Value = 1300
$ws.Cells.Item(100, 11).Value = 7483245.5
$ws.Cells.Item(100, 12).Value = 1300
$ws.Cells.Item(100, 13).Value = -7482704.5
$ws.Cells.Item(100, 14).Value = -2382

$ws.Cells.Item(108, 8).Value = 43800
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 43800
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = 43800
$ws.Cells.Item(108, 14).Value = -51480

$ws.Cells.Item(113, 8).Value = 1329.375
$ws.Cells.Item(113, 9).Value = 1351.5385
$ws.Cells.Item(113, 10).Value = 1233.3334
$ws.Cells.Item(113, 11).Value = 1351.5385
$ws.Cells.Item(113, 12).Value = 1233.3334
$ws.Cells.Item(113, 13).Value = 818.4614999999999
$ws.Cells.Item(113, 14).Value = -5573.3334
